$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("detection_template_csv")

# Update the formula in O3 to use a VLOOKUP against the 'lookup' sheet for the assay name,
# and drop the trailing "from an organism" suffix.
$ws.Range("O3").Formula = '=IF($D3="count","a count of the number of ",IF($D3="boolean","a categorical measurement datum","a data item")&" that is about ")&$H3&" and is the specified output of some "&IF(ISNA(VLOOKUP(C3,lookup!A2:B4,2,FALSE)=TRUE),C3,VLOOKUP(C3,lookup!A2:B4,2))&", which achieves an organism identification objective and has as specified input a "&$B3&" specimen"'

# Update the formula in P3 similarly.
$ws.Range("P3").Formula = '="("&IF($D3="count","count and",IF($D3="boolean","''categorical measurement datum'' and","''data item'' and")&" ''is about'' some ")&"''"&$H3&"'') and is_specified_output_of some ((''"&IF(ISNA(VLOOKUP(C3,lookup!A2:B4,2,FALSE)=TRUE),C3,VLOOKUP(C3,lookup!A2:B4,2))&"'' and achieves_planned_objective some ''organism identification objective'') and has_specified_input some ''"&$B3&" specimen'')"'

# Move the active cell selection from O3 to O4.
$ws.Activate()
$ws.Range("O4").Select()
